$d = $word.ActiveDocument

$d.Content.Find.Execute("74×32=2368", $true, $false, $false, $false, $false, $true, 1, $false, "81×43=3483", 2) | Out-Null
$d.Content.Find.Execute("75×19=1425", $true, $false, $false, $false, $false, $true, 1, $false, "61×73=4453", 2) | Out-Null
$d.Content.Find.Execute("91×21=1911", $true, $false, $false, $false, $false, $true, 1, $false, "47×89=4183", 2) | Out-Null
$d.Content.Find.Execute("91×96=8736", $true, $false, $false, $false, $false, $true, 1, $false, "27×64=1728", 2) | Out-Null
$d.Content.Find.Execute("79×97=7663", $true, $false, $false, $false, $false, $true, 1, $false, "54×58=3132", 2) | Out-Null
$d.Content.Find.Execute("52×24=1248", $true, $false, $false, $false, $false, $true, 1, $false, "94×90=8460", 2) | Out-Null
$d.Content.Find.Execute("96×85=8160", $true, $false, $false, $false, $false, $true, 1, $false, "46×33=1518", 2) | Out-Null
$d.Content.Find.Execute("58×96=5568", $true, $false, $false, $false, $false, $true, 1, $false, "94×44=4136", 2) | Out-Null
$d.Content.Find.Execute("66×68=4488", $true, $false, $false, $false, $false, $true, 1, $false, "72×51=3672", 2) | Out-Null
$d.Content.Find.Execute("33×99=3267", $true, $false, $false, $false, $false, $true, 1, $false, "23×13=299", 2) | Out-Null
$d.Content.Find.Execute("72×92=6624", $true, $false, $false, $false, $false, $true, 1, $false, "41×89=3649", 2) | Out-Null
$d.Content.Find.Execute("48×92=4416", $true, $false, $false, $false, $false, $true, 1, $false, "90×67=6030", 2) | Out-Null
$d.Content.Find.Execute("42×81=3402", $true, $false, $false, $false, $false, $true, 1, $false, "39×46=1794", 2) | Out-Null
$d.Content.Find.Execute("79×90=7110", $true, $false, $false, $false, $false, $true, 1, $false, "73×92=6716", 2) | Out-Null
$d.Content.Find.Execute("38×14=532", $true, $false, $false, $false, $false, $true, 1, $false, "59×95=5605", 2) | Out-Null
$d.Content.Find.Execute("56×39=2184", $true, $false, $false, $false, $false, $true, 1, $false, "70×46=3220", 2) | Out-Null
$d.Content.Find.Execute("80×71=5680", $true, $false, $false, $false, $false, $true, 1, $false, "65×51=3315", 2) | Out-Null
$d.Content.Find.Execute("75×57=4275", $true, $false, $false, $false, $false, $true, 1, $false, "29×41=1189", 2) | Out-Null
$d.Content.Find.Execute("43×32=1376", $true, $false, $false, $false, $false, $true, 1, $false, "15×46=690", 2) | Out-Null
$d.Content.Find.Execute("87×15=1305", $true, $false, $false, $false, $false, $true, 1, $false, "40×51=2040", 2) | Out-Null
$d.Content.Find.Execute("20×99=1980", $true, $false, $false, $false, $false, $true, 1, $false, "83×85=7055", 2) | Out-Null
$d.Content.Find.Execute("62×44=2728", $true, $false, $false, $false, $false, $true, 1, $false, "83×49=4067", 2) | Out-Null
$d.Content.Find.Execute("34×91=3094", $true, $false, $false, $false, $false, $true, 1, $false, "69×11=759", 2) | Out-Null
$d.Content.Find.Execute("45×11=495", $true, $false, $false, $false, $false, $true, 1, $false, "38×63=2394", 2) | Out-Null
$d.Content.Find.Execute("16×52=832", $true, $false, $false, $false, $false, $true, 1, $false, "77×97=7469", 2) | Out-Null
